# edit.ps1
# Applies the RECO_holdings.xlsx update:
#  - Refresh the "as of" date in the confidential disclaimer (A41)
#  - Update Weight (col D) and Percent Change (col E) figures for rows 2-38
#
# The sheet is protected (no password known to us), so we unprotect it,
# make the edits, and re-protect it to restore the sheet's protected state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Update confidential disclaimer text (as-of date 2021-05-03 -> 2021-05-04) ---
$ws.Range("A41").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."
# Re-setting a multi-line value can mark the row as a custom (auto-fit) height;
# AutoFit puts it back to the sheet's natural/default height so row 41 stays unchanged.
$ws.Rows(41).AutoFit()

# --- Update Weight / Percent Change figures ---
$ws.Range("D2").Value = 0.03118664532805398
$ws.Range("E2").Value = -0.006801534705266787
$ws.Range("D3").Value = 0.02933858523754682
$ws.Range("E3").Value = -0.0253446152623843
$ws.Range("D4").Value = 0.02899030147232931
$ws.Range("E4").Value = 0.0007638446849140834
$ws.Range("D5").Value = 0.06578134345072095
$ws.Range("E5").Value = -0.02203461401037654
$ws.Range("D6").Value = 0.01570365460972891
$ws.Range("E6").Value = -0.0001484340210777813
$ws.Range("D7").Value = 0.01601716827233965
$ws.Range("E7").Value = -0.02408498872153098
$ws.Range("D8").Value = 0.02859850651725753
$ws.Range("E8").Value = 0.008069117287472505
$ws.Range("D9").Value = 0.03437383998487987
$ws.Range("E9").Value = -0.001695298372513698
$ws.Range("D10").Value = 0.02966783285720307
$ws.Range("E10").Value = -0.01164777749405821
$ws.Range("D11").Value = 0.03137020816025172
$ws.Range("E11").Value = -0.0121797928134888
$ws.Range("D12").Value = 0.01253666157676394
$ws.Range("E12").Value = -0.04090486519987613
$ws.Range("D13").Value = 0.01396709194322334
$ws.Range("E13").Value = -0.0128504672897195
$ws.Range("D14").Value = 0.01529146378441293
$ws.Range("E14").Value = -0.02258580828739087
$ws.Range("D15").Value = 0.008907362151597259
$ws.Range("E15").Value = 0.00680390788555485
$ws.Range("D16").Value = 0.007636213819425843
$ws.Range("E16").Value = 0.02083333333333348
$ws.Range("D17").Value = 0.03042519950560411
$ws.Range("E17").Value = -0.0294065069717554
$ws.Range("D18").Value = 0.02636700406615955
$ws.Range("E18").Value = -0.03403565640194495
$ws.Range("D19").Value = 0.03118237190762186
$ws.Range("E19").Value = -0.01709960754999074
$ws.Range("D20").Value = 0.03132999915891316
$ws.Range("E20").Value = -0.01308202616405219
$ws.Range("D21").Value = 0.04652531098360348
$ws.Range("E21").Value = -0.01708438231941789
$ws.Range("D22").Value = 0.03400865678431718
$ws.Range("E22").Value = -0.0009138679460819477
$ws.Range("D23").Value = 0.03207687650261718
$ws.Range("E23").Value = 0.007569564295879072
$ws.Range("D24").Value = 0.03043646579583423
$ws.Range("E24").Value = 0.01418724870763932
$ws.Range("D25").Value = 0.01536857959857431
$ws.Range("E25").Value = -0.04459105903765215
$ws.Range("D26").Value = 0.01501835336952403
$ws.Range("E26").Value = -0.03499922396399213
$ws.Range("D27").Value = 0.03092888037744404
$ws.Range("E27").Value = 0.006908462867012188
$ws.Range("D28").Value = 0.02985489212429982
$ws.Range("E28").Value = 0.00806787424526334
$ws.Range("D29").Value = 0.0293537364554425
$ws.Range("E29").Value = -0.0161597713015168
$ws.Range("D30").Value = 0.02842465600422369
$ws.Range("E30").Value = -0.01112531008043283
$ws.Range("D31").Value = 0.03458382032520341
$ws.Range("E31").Value = -0.03272279980453952
$ws.Range("D32").Value = 0.0311510982399141
$ws.Range("E32").Value = -0.01255230125522988
$ws.Range("D33").Value = 0.02927720338043096
$ws.Range("E33").Value = -0.005002587545282045
$ws.Range("D34").Value = 0.03153784278902065
$ws.Range("E34").Value = 0.01340231584134033
$ws.Range("D35").Value = 0.03016024743881288
$ws.Range("E35").Value = 0
$ws.Range("D36").Value = 0.0297896253395184
$ws.Range("E36").Value = 0.01382368283776714
$ws.Range("D37").Value = 0.03283230068718544
$ws.Range("E37").Value = 0.001561908368042486
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = -0.009361121702934683

$ws.Protect()
